# The commit inserts one new weekly observation row for "Orégano" at
# Mercado Mayorista Lo Valledor de Santiago. The new record belongs
# chronologically between the existing rows that are currently at
# sheet rows 60 and 61, so a new row is inserted at row 61, pushing
# all subsequent rows (old 61..129) down by one (new 62..130).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 61 (shifts rows 61-129 down to 62-130)
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new observation
$ws.Cells.Item(61, 1).Value  = 6
$ws.Cells.Item(61, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(61, 3).Value  = "Metropolitana"
$ws.Cells.Item(61, 4).Value  = 44539
$ws.Cells.Item(61, 5).Value  = 13
$ws.Cells.Item(61, 6).Value  = 100112029
$ws.Cells.Item(61, 7).Value  = "Orégano"
$ws.Cells.Item(61, 8).Value  = "Sin especificar"
$ws.Cells.Item(61, 9).Value  = "Primera"
$ws.Cells.Item(61, 10).Value = 35
$ws.Cells.Item(61, 11).Value = 8000
$ws.Cells.Item(61, 12).Value = 9000
$ws.Cells.Item(61, 13).Value = 8457
$ws.Cells.Item(61, 14).Value = "`$/docena de atados"
$ws.Cells.Item(61, 15).Value = "Región Metropolitana"
$ws.Cells.Item(61, 16).Value = 2819
$ws.Cells.Item(61, 17).Value = 3
$ws.Cells.Item(61, 18).Value = "Hortaliza"

# Give the new date cell the same date number format used by the rest
# of column D (style index 2 in the original workbook -> "YYYY-MM-DD HH:MM:SS")
$ws.Cells.Item(61, 4).NumberFormat = $ws.Cells.Item(62, 4).NumberFormat
